$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")
$ws.Hyperlinks.Delete() | Out-Null

$data = @(
    @('2','2025-10-22 18:28:54','【短期/長期OK】AIエンジニアパートナー募集!業務効率化・自動化支援','システム開発','100,000 円 ~ 200,000 円 / 固定','期限情報なし','https://www.lancers.jp/work/detail/5418449','438','🔥AI,Ai ◆効率化,自動化'),
    @('3','2025-10-22 18:28:54','【ECシステム開発】販売データ分析・AI提案・競合監視を統合した販売支援システム構築','システム開発','500,000 円 ~ 1,000,000 円 / 固定','期限情報なし','https://www.lancers.jp/work/detail/5418284','410','🔥AI,Ai ◆開発,システム開発'),
    @('4','2025-10-22 18:28:54','外国人雇用事業所向けAI日報アプリ開発の依頼','システム開発','200,000 円 ~ 300,000 円 / 固定','期限情報なし','https://www.lancers.jp/work/detail/5418427','378','🔥AI,Ai ◆開発 ◇アプリ'),
    @('5','2025-10-22 18:28:54','【AIで開発生産性を革新】AI活用推進エンジニア募集(副業・業務委託)','システム開発','500,000 円 ~ 1,000,000 円 / 固定','期限情報なし','https://www.lancers.jp/work/detail/5417964','375','🔥AI,Ai ◆開発'),
    @('6','2025-10-22 18:28:54','【業務委託/副業可】AI SaaS開発を牽引するCTO候補を募集','システム開発','500,000 円 ~ 1,000,000 円 / 固定','期限情報なし','https://www.lancers.jp/work/detail/5417967','375','🔥AI,Ai ◆開発'),
    @('7','2025-10-22 18:28:54','【継続依頼あり】AI×業務効率化のスペシャリスト募集!','システム開発','20,000 円 ~ 50,000 円 / 固定','期限情報なし','https://www.lancers.jp/work/detail/5418075','373','🔥AI,Ai ◆効率化'),
    @('8','2025-10-22 18:28:54','【AI技術顧問/戦略アドバイザー募集】最先端AIで事業の非連続な成長を牽引するエキスパート求む','システム開発','500,000 円 ~ 1,000,000 円 / 固定','期限情報なし','https://www.lancers.jp/work/detail/5417960','310','🔥AI,Ai'),
    @('9','2025-10-22 18:28:54','GASと生成AIを活用したスプレッドシートの作り方レクチャー','システム開発','20,000 円 ~ 50,000 円 / 固定','期限情報なし','https://www.lancers.jp/work/detail/5418291','298','🔥AI,Ai'),
    @('10','2025-10-22 18:28:54','【せどり×ツール製作】APIを使用したせどりツールを製作できるエンジニアさんを募集します♪','システム開発','20,000 円 ~ 50,000 円 / 固定','期限情報なし','https://www.lancers.jp/work/detail/5217096','243','🔥API ◆ツール'),
    @('11','2025-10-22 18:28:54','【スプレッドシート自動化】業務効率化の専門家を募集します','システム開発','20,000 円 ~ 50,000 円 / 固定','期限情報なし','https://www.lancers.jp/work/detail/5418424','148','◆効率化,自動化'),
    @('12','2025-10-22 18:28:54','【急募】キントーン見積書をエクセルに変換するツール開発','システム開発','10,000 円 ~ 20,000 円 / 固定','期限情報なし','https://www.lancers.jp/work/detail/5418067','120','◆ツール,開発'),
    @('13','2025-10-22 18:28:54','【急募】Webアプリ開発エンジニア募集!フルリモート可','システム開発','200,000 円 ~ 300,000 円 / 固定','期限情報なし','https://www.lancers.jp/work/detail/5411585','93','◆開発 ◇アプリ'),
    @('14','2025-10-22 18:28:54','大手クレジットカード企業向け、Google Cloudを利用したアジャイル開発共通基盤案件_ワーカー','システム開発','500,000 円 ~ 1,000,000 円 / 固定','期限情報なし','https://www.lancers.jp/work/detail/5418318','75','◆開発'),
    @('15','2025-10-22 18:28:54','大手クレジットカード企業向け、Google Cloudを利用したアジャイル開発共通基盤案件','システム開発','500,000 円 ~ 1,000,000 円 / 固定','期限情報なし','https://www.lancers.jp/work/detail/5418320','75','◆開発'),
    @('16','2025-10-22 18:28:54','IISIA公式サイト WordPressアップデート&AWS運用整備 実施要領書','システム開発','300,000 円 ~ 500,000 円 / 固定','期限情報なし','https://www.lancers.jp/work/detail/5418421','65','◇サイト ○WordPress'),
    @('17','2025-10-22 18:28:54','【急募】WEB会計アプリのテストユーザーを募集します!','システム開発','5,000 円 ~ 10,000 円 / 固定','期限情報なし','https://www.lancers.jp/work/detail/5418565','30','◇アプリ'),
    @('18','2025-10-22 18:28:54','進行管理およびチームディレクションを担当','システム開発','~ 5,000 円 / 固定','期限情報なし','https://www.lancers.jp/work/detail/5418064','30','◇管理'),
    @('19','2025-10-22 18:28:54','自社HPに見積自動受付システムを設置したい','システム開発','200,000 円 ~ 300,000 円 / 固定','期限情報なし','https://www.lancers.jp/work/detail/5418456','33',''),
    @('20','2025-10-22 18:28:54','【WP安全アップデート+AWS運用まで一括/haradatakeo.com(~45万円以下)】','システム開発','300,000 円 ~ 500,000 円 / 固定','期限情報なし','https://www.lancers.jp/work/detail/5418426','25',''),
    @('21','2025-10-22 18:28:54','サブスクペイからCSVデータをダウンロードし、データベース同期するプログラムの作成','システム開発','50,000 円 ~ 100,000 円 / 固定','期限情報なし','https://www.lancers.jp/work/detail/5418241','18',''),
    @('22','2025-10-22 18:28:54','ハードウェアの設定設置と保守サポート依頼|東京周辺対応可能な方','システム開発','100,000 円 ~ 200,000 円 / 固定','期限情報なし','https://www.lancers.jp/work/detail/5418084','18',''),
    @('23','2025-10-22 18:28:54','google workspace の設定方法を教えて下さい。','システム開発','~ 5,000 円 / 固定','期限情報なし','https://www.lancers.jp/work/detail/5418606','10',''),
    @('24','2025-10-22 18:28:54','【急募】Meta広告のコンバージョン計測設定をお手伝いください!','システム開発','5,000 円 ~ 10,000 円 / 固定','期限情報なし','https://www.lancers.jp/work/detail/5418533','10','')
)

foreach ($row in $data) {
    $r = [int]$row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Hyperlinks.Add($ws.Cells.Item($r, 6), $row[6]) | Out-Null
    $ws.Cells.Item($r, 6).Value = $row[6]
    $ws.Cells.Item($r, 6).Style = "Hyperlink"
    $ws.Cells.Item($r, 7).Value = [double]$row[7]
    if ($row[8] -ne "") {
        $ws.Cells.Item($r, 8).Value = $row[8]
    } else {
        $ws.Cells.Item($r, 8).Value = ""
    }
}
